$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

$ws.Range("B31").Value = "Yes"

$ws.Range("A46").Value = "SRO Schedule Confidence"
$ws.Range("B46").Value = "Green"
$ws.Range("C46").Value = "Red"
$ws.Range("D46").Value = "Amber"
$ws.Range("E46").Value = "Amber/Red"
$ws.Range("F46").Value = "Amber"
$ws.Range("G46").Value = "Amber/Green"



$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("F37").Select()
